$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (issue number + date range) ---
$ws.Range("A8").Value = "Volume 32   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/10/2025  Through  2/16/2025"

# --- Cells changing between numeric and N/A-text representation ---
# For these we first fix up the value/type, then re-apply the donor cell's
# number format (via PasteSpecial formats-only) so the style index matches
# the style used by other cells of the same kind in the sheet.

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4122)

$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4122)

$ws.Range("D15").Value = 1
$ws.Range("F15").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").Value = -100
$ws.Range("H15").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("C18").Value = 3
$ws.Range("F18").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("D18").Value = 1
$ws.Range("F18").Copy()
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("E18").Value = 200
$ws.Range("H18").Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D27").Value = 1
$ws.Range("F27").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = -100
$ws.Range("H27").Copy()
$ws.Range("E27").PasteSpecial(-4122)

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 50
$ws.Range("M15").Value = -25
$ws.Range("N15").Value = -50
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = -62.5
$ws.Range("I16").Value = 25
$ws.Range("J16").Value = 29
$ws.Range("K16").Value = -13.793103448275
$ws.Range("L16").Value = 13.636363636363
$ws.Range("M16").Value = -41.860465116279
$ws.Range("N16").Value = -81.481481481481
$ws.Range("C17").Value = 11
$ws.Range("E17").Value = 120
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = 3.846153846153
$ws.Range("I17").Value = 43
$ws.Range("J17").Value = 36
$ws.Range("K17").Value = 19.444444444444
$ws.Range("L17").Value = -6.521739130434
$ws.Range("M17").Value = 152.941176470588
$ws.Range("N17").Value = 4.878048780487
$ws.Range("F18").Value = 7
$ws.Range("H18").Value = 133.333333333333
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 16
$ws.Range("K18").Value = -12.5
$ws.Range("L18").Value = -17.647058823529
$ws.Range("M18").Value = -50
$ws.Range("N18").Value = -92.090395480226
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -40
$ws.Range("F19").Value = 29
$ws.Range("H19").Value = -40.816326530612
$ws.Range("I19").Value = 53
$ws.Range("J19").Value = 67
$ws.Range("K19").Value = -20.895522388059
$ws.Range("L19").Value = -38.372093023255
$ws.Range("M19").Value = 43.243243243243
$ws.Range("N19").Value = -27.397260273972
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 29
$ws.Range("J20").Value = 53
$ws.Range("K20").Value = -45.283018867924
$ws.Range("L20").Value = -35.555555555555
$ws.Range("M20").Value = -25.641025641025
$ws.Range("N20").Value = -93.933054393305
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 27.272727272727
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 131
$ws.Range("H21").Value = -30.534351145038
$ws.Range("I21").Value = 167
$ws.Range("J21").Value = 204
$ws.Range("K21").Value = -18.137254901960
$ws.Range("L21").Value = -23.394495412844
$ws.Range("M21").Value = -0.595238095238
$ws.Range("N21").Value = -81.788440567066
$ws.Range("G22").Value = 1
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 47.619047619047
$ws.Range("F24").Value = 105
$ws.Range("G24").Value = 92
$ws.Range("H24").Value = 14.130434782608
$ws.Range("I24").Value = 171
$ws.Range("J24").Value = 143
$ws.Range("K24").Value = 19.580419580419
$ws.Range("L24").Value = 14.765100671140
$ws.Range("M24").Value = 122.077922077922
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 144.444444444444
$ws.Range("F25").Value = 58
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = 16
$ws.Range("I25").Value = 92
$ws.Range("J25").Value = 65
$ws.Range("K25").Value = 41.538461538461
$ws.Range("L25").Value = 50.819672131147
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -16.666666666666
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 34
$ws.Range("H26").Value = 2.941176470588
$ws.Range("I26").Value = 66
$ws.Range("J26").Value = 61
$ws.Range("K26").Value = 8.196721311475
$ws.Range("L26").Value = 3.125
$ws.Range("M26").Value = 10
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 2
$ws.Range("K27").Value = 100
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -28.571428571428
$ws.Range("I28").Value = 7
$ws.Range("J28").Value = 7
$ws.Range("L28").Value = -46.153846153846
